$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 359, shifting existing rows 359:381 down to 361:383.
$ws.Rows.Item(359).Resize(2).Insert()

# New row 359: Zafiro rojo
$ws.Cells.Item(359, 1).Value = 11
$ws.Cells.Item(359, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(359, 3).Value = "Bíobío"
$ws.Cells.Item(359, 4).Value = 44826
$ws.Cells.Item(359, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(359, 5).Value = 8
$ws.Cells.Item(359, 6).Value = 100112002
$ws.Cells.Item(359, 7).Value = "Pimiento"
$ws.Cells.Item(359, 8).Value = "Zafiro rojo"
$ws.Cells.Item(359, 9).Value = "Primera"
$ws.Cells.Item(359, 10).Value = 100
$ws.Cells.Item(359, 11).Value = 16000
$ws.Cells.Item(359, 12).Value = 17000
$ws.Cells.Item(359, 13).Value = 16500
$ws.Cells.Item(359, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(359, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(359, 16).Value = 1100
$ws.Cells.Item(359, 17).Value = 15
$ws.Cells.Item(359, 18).Value = "Hortaliza"

# New row 360: Zafiro verde
$ws.Cells.Item(360, 1).Value = 11
$ws.Cells.Item(360, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(360, 3).Value = "Bíobío"
$ws.Cells.Item(360, 4).Value = 44826
$ws.Cells.Item(360, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(360, 5).Value = 8
$ws.Cells.Item(360, 6).Value = 100112002
$ws.Cells.Item(360, 7).Value = "Pimiento"
$ws.Cells.Item(360, 8).Value = "Zafiro verde"
$ws.Cells.Item(360, 9).Value = "Primera"
$ws.Cells.Item(360, 10).Value = 100
$ws.Cells.Item(360, 11).Value = 16000
$ws.Cells.Item(360, 12).Value = 17000
$ws.Cells.Item(360, 13).Value = 16500
$ws.Cells.Item(360, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(360, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(360, 16).Value = 1100
$ws.Cells.Item(360, 17).Value = 15
$ws.Cells.Item(360, 18).Value = "Hortaliza"
